$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was updated with a new weekly price record that lands at
# row 40 (immediately after the existing row 39). Insert a new row there,
# which shifts every following row down by one (old row 40 -> new row 41,
# ..., old row 157 -> new row 158), and populate the new row with the
# latest week's values.
$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "Terminal La Palmera de La Serena"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44972
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 100112052
$ws.Range("G40").Value = "Albahaca"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 840
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = 5500
$ws.Range("N40").Value = "`$/docena de matas"
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 917
$ws.Range("Q40").Value = 6
$ws.Range("R40").Value = "Hortaliza"
